$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 ("مبرد قدم" line item): update the transaction-count fields (H, Q)
# and the selling price (P). These are stored as text in the sheet, so for
# the numeric-looking "60.0000" we round-trip the number format to keep it
# a text value instead of letting Excel coerce it to a number.
$ws.Range("H17").Value = "1:0"

$ws.Range("P17").NumberFormat = "@"
$ws.Range("P17").Value = "60.0000"
$ws.Range("P17").NumberFormat = "0.00"

$ws.Range("Q17").Value = "4:0"

# Row 18's running total increases by the same amount the price went up.
$ws.Range("P18").Value = 1077.3199999999999
